$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("loginTest")
$wsData  = $wb.Worksheets.Item("Sheet1")

# --- Rework the "Sheet1" tab ---
# Remove the existing hyperlinks (and their styled cells) then clear all data.
$wsData.Hyperlinks.Delete() | Out-Null
$wsData.UsedRange.Delete() | Out-Null

# Write the new, simplified A1:C5 table (username/password + pass/fail status).
$rows = @(
  @("username",         "password",         "status"),
  @("wrong_username",   "wrong_password",   "failed"),
  @("wrong_username",   "wrong_password",   "failed"),
  @("wrong_username",   "wrong_password",   "failed"),
  @("correct_username", "correct_password", "passed")
)
for ($r = 0; $r -lt $rows.Length; $r++) {
  for ($c = 0; $c -lt 3; $c++) {
    $wsData.Cells.Item($r + 1, $c + 1).Value = $rows[$r][$c]
  }
}

# Approximate column widths for the new A/B columns.
$wsData.Columns.Item(1).ColumnWidth = 16.8
$wsData.Columns.Item(2).ColumnWidth = 16.65

# Remove the now-unused "Hyperlink" cell style left over from the deleted hyperlinks.
$wb.Styles.Item("Hyperlink").Delete() | Out-Null

# --- Selection / active sheet bookkeeping ---
$wsLogin.Range("A1:B5").Select() | Out-Null
$wsData.Range("C5").Select() | Out-Null
$wsData.Activate() | Out-Null
